# Insert a new price-record row at row 173 ("Terminal La Palmera de La Serena - Poroto verde"
# subset), which pushes the existing rows 173-248 down to 174-249 and grows the used range
# from A1:R248 to A1:R249.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 173..248 down by one row.
$ws.Rows("173:173").Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(173, 1).Value = 8
$ws.Cells.Item(173, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(173, 3).Value = "Coquimbo"
$ws.Cells.Item(173, 4).Value = 44755
$ws.Cells.Item(173, 5).Value = 4
$ws.Cells.Item(173, 6).Value = 100112031
$ws.Cells.Item(173, 7).Value = "Poroto verde"
$ws.Cells.Item(173, 8).Value = "Magnum"
$ws.Cells.Item(173, 9).Value = "Primera"
$ws.Cells.Item(173, 10).Value = 480
$ws.Cells.Item(173, 11).Value = 34000
$ws.Cells.Item(173, 12).Value = 35000
$ws.Cells.Item(173, 13).Value = 34500
$ws.Cells.Item(173, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(173, 15).Value = "Perú"
$ws.Cells.Item(173, 16).Value = 1380
$ws.Cells.Item(173, 17).Value = 25
$ws.Cells.Item(173, 18).Value = "Hortaliza"
